$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: Leftfield Brewery (new data replacing the previously-empty placeholder row)
$ws.Range("A16").Value = "Leftfield Brewery"
$ws.Range("B16").Value = "Leslieville"
$ws.Range("C16").Value = "Fenway Red - really good beer"
# Leading apostrophe restores the original cell's "quote prefix" formatting (s=2)
# while still storing the clean text value "Brewery".
$ws.Range("D16").Value = "'Brewery"
$ws.Range("E16").Value = 43.674216711648498
$ws.Range("F16").Value = -79.330432202501399

# Row 17: La Banane
$ws.Range("A17").Value = "La Banane"
$ws.Range("B17").Value = "Trinity Bellwoods"
$ws.Range("C17").Value = "Raw Bar, Maitake Mushroom was a sneaky star, amazing cocktails"
$ws.Range("D17").Value = "French"
$ws.Range("E17").Value = 43.649649350878001
$ws.Range("F17").Value = -79.4204187426581

# Row 18: Lamoon Thai Café
$ws.Range("A18").Value = "Lamoon Thai Café"
$ws.Range("B18").Value = "Leslieville"
$ws.Range("C18").Value = "Cuteness levels unmatched, lunch specials are great, thai tea c'mon now"
$ws.Range("D18").Value = "Thai"
$ws.Range("E18").Value = 43.666404968196701
$ws.Range("F18").Value = -79.348041345209097

# Column B widened to fit the new longer location text ("Trinity Bellwoods")
$ws.Columns("B").ColumnWidth = 14.5

# Update the active selection to mirror where the author's cursor ended up
$ws.Range("C19").Select() | Out-Null
